$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 cleanup: remove the now-empty placeholder cells (D3, G3, I3, K3) ---
$ws.Range("D3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("K3").ClearContents()

# --- Row 4: Bilbo Baggins ---
$ws.Range("A4").Value = "Bilbo"
$ws.Range("B4").Value = "Baggins"
$ws.Range("C4").Value = "Party Planning"
$ws.Range("D4").Value = "n/a"
$ws.Range("E4").Value = "Mithrandir"
# Leading apostrophe forces these to stay plain text instead of being
# auto-parsed into Excel date serials; resetting the style back to Normal
# afterwards drops the quote-prefix formatting so the cell is plain text.
$ws.Range("F4").Value = "'04/24/2022"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'10/31/2026"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "bbaggins@outlook.com"
$ws.Range("I4").Value = "n/a"
$ws.Range("J4").Value = "555-456-8764"
$ws.Range("K4").Value = "n/a"

# --- Row 5: Peregrin Took ---
$ws.Range("A5").Value = "Peregrin"
$ws.Range("B5").Value = "Took"
$ws.Range("C5").Value = "Farming Consulting"
$ws.Range("E5").Value = "Meriadoc Brandybuck"
$ws.Range("F5").Value = "'07/16/2023"
$ws.Range("F5").Style = "Normal"
$ws.Range("H5").Value = "ptook@outlook.com"
$ws.Range("I5").Value = 15553490293
$ws.Range("J5").Value = "555-777-7676"
$ws.Range("K5").Value = "n/a"
